$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated prediction/error/cross-entropy-loss/success% results with new
# toy dataset values (NCDE/NODE architecture testing datasets).

$ws.Range("D2").Value = 0.8142413394034955
$ws.Range("E2").Value = 0.8142413394034955
$ws.Range("D3").Value = 0.03929915539604146
$ws.Range("E3").Value = 0.03929915539604146
$ws.Range("D4").Value = 0.04385414019569965
$ws.Range("E4").Value = 0.04385414019569965
$ws.Range("D5").Value = 0.00209750943385401
$ws.Range("E5").Value = 0.00209750943385401
$ws.Range("D6").Value = 0.9465492369353946
$ws.Range("E6").Value = 0.9465492369353946
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("D8").Value = 0.9955583887171143
$ws.Range("E8").Value = 0.004441611282885671
$ws.Range("D9").Value = 0.9999999187017909
$ws.Range("E9").Value = [double]"8.129820905367069E-08"
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 0.02983007302867908
$ws.Range("E10").Value = 0.9701699269713209
$ws.Range("D11").Value = 0.9472138594349139
$ws.Range("E11").Value = 0.05278614056508613
$ws.Range("F11").Value = 0.8270258903503418
$ws.Range("G11").Value = 0.7
$ws.Range("D12").Value = 0.8680562665205624
$ws.Range("E12").Value = 0.8680562665205624
$ws.Range("D13").Value = 0.005964977039969782
$ws.Range("E13").Value = 0.005964977039969782
$ws.Range("D14").Value = 0.0044990834001208
$ws.Range("E14").Value = 0.0044990834001208
$ws.Range("D15").Value = [double]"9.278267216329316E-05"
$ws.Range("E15").Value = [double]"9.278267216329316E-05"
$ws.Range("D16").Value = 0.964949302692886
$ws.Range("E16").Value = 0.964949302692886
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("D18").Value = 0.9868002054308599
$ws.Range("E18").Value = 0.01319979456914011
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = [double]"2.938384125245676E-07"
$ws.Range("E19").Value = 0.9999997061615875
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = 0.01457421051713957
$ws.Range("E20").Value = 0.9854257894828604
$ws.Range("D21").Value = 0.9718275914417014
$ws.Range("E21").Value = 0.02817240855829861
$ws.Range("F21").Value = 2.469752550125122
$ws.Range("G21").Value = 0.6
